$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 11.32
$ws.Range("D2").Value = 11.42

# Row 3
$ws.Range("B3").Value = 8.68
$ws.Range("D3").Value = 10.37
$ws.Range("E3").Value = 10.82

# Row 4
$ws.Range("B4").Value = 8.58
$ws.Range("C4").Value = 9.630000000000001
$ws.Range("E4").Value = 10.6
$ws.Range("F4").Value = 10.1

# Row 5
$ws.Range("C5").Value = 9.18
$ws.Range("D5").Value = 9.4
$ws.Range("F5").Value = 10.31

# Row 6
$ws.Range("D6").Value = 9.9
$ws.Range("E6").Value = 9.69
$ws.Range("G6").Value = 10.4
$ws.Range("H6").Value = 10.89
$ws.Range("I6").Value = 8.76

# Row 7
$ws.Range("F7").Value = 9.6
$ws.Range("H7").Value = 9.52
$ws.Range("J7").Value = 9.58

# Row 8
$ws.Range("F8").Value = 9.109999999999999
$ws.Range("G8").Value = 10.48

# Row 9
$ws.Range("F9").Value = 11.24

# Row 10
$ws.Range("G10").Value = 10.42
